$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 837.1
$ws.Range("J18").Value = 2004
$ws.Range("L18").Value = 2004
$ws.Range("N18").Value = -2572
$ws.Range("I40").Value = 1979.8
$ws.Range("K40").Value = 1979.8
$ws.Range("M40").Value = -1804.8
$ws.Range("H111").Value = 1622.3636
$ws.Range("I111").Value = 1049.625
$ws.Range("K111").Value = 3148.875
$ws.Range("M111").Value = -81.875
$ws.Range("H135").Value = 1415.6
$ws.Range("I135").Value = 1415.6
$ws.Range("K135").Value = 12740.4
$ws.Range("M135").Value = -10205.4
$ws.Range("H137").Value = 2550.1304
$ws.Range("I137").Value = 2877.1875
$ws.Range("K137").Value = 8631.5625
$ws.Range("M137").Value = -6081.5625
$ws.Range("H138").Value = 24394024
$ws.Range("J138").Value = 3961.125
$ws.Range("L138").Value = 11883.375
$ws.Range("N138").Value = -22163.375
$ws.Range("H141").Value = 3351.9167
$ws.Range("I141").Value = 1967.2858
$ws.Range("K141").Value = 5901.857400000001
$ws.Range("M141").Value = -721.8574000000008

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2558.1875
$ws.Range("I61").Value = 1449
$ws.Range("J61").Value = 4998.4
$ws.Range("K61").Value = 1449
$ws.Range("L61").Value = 4998.4
$ws.Range("M61").Value = -1237
$ws.Range("N61").Value = -5422.4
$ws.Range("H74").Value = 2817.3096
$ws.Range("I74").Value = 2322.3794
$ws.Range("K74").Value = 2322.3794
$ws.Range("M74").Value = -1448.3794
$ws.Range("H77").Value = 2817.3096
$ws.Range("I77").Value = 2322.3794
$ws.Range("K77").Value = 11611.897
$ws.Range("M77").Value = -7243.896999999999
$ws.Range("H82").Value = 60000
$ws.Range("J82").Value = 60000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60722
$ws.Range("H85").Value = 60000
$ws.Range("J85").Value = 60000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62496
$ws.Range("H122").Value = 3875.9443
$ws.Range("I122").Value = 3455.5186
$ws.Range("K122").Value = 10366.5558
$ws.Range("M122").Value = -7916.5558
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2558.1875
$ws.Range("I136").Value = 1449
$ws.Range("J136").Value = 4998.4
$ws.Range("K136").Value = 4347
$ws.Range("L136").Value = 14995.2
$ws.Range("M136").Value = -1797
$ws.Range("N136").Value = -20095.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3062.9473
$ws.Range("I86").Value = 2899.818
$ws.Range("J86").Value = 3287.25
$ws.Range("K86").Value = 2899.818
$ws.Range("L86").Value = 3287.25
$ws.Range("M86").Value = -1776.818
$ws.Range("N86").Value = -5533.25
$ws.Range("H89").Value = 3062.9473
$ws.Range("I89").Value = 2899.818
$ws.Range("J89").Value = 3287.25
$ws.Range("K89").Value = 14499.09
$ws.Range("L89").Value = 16436.25
$ws.Range("M89").Value = -8883.09
$ws.Range("N89").Value = -27668.25
$ws.Range("H134").Value = 3801.4897
$ws.Range("I134").Value = 3537.7026
$ws.Range("K134").Value = 10613.1078
$ws.Range("M134").Value = -8078.1078

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1696.75
$ws.Range("I31").Value = 1349.1
$ws.Range("K31").Value = 1349.1
$ws.Range("M31").Value = -1054.1
$ws.Range("H34").Value = 1696.75
$ws.Range("I34").Value = 1349.1
$ws.Range("K34").Value = 1349.1
$ws.Range("M34").Value = -1147.1
$ws.Range("H86").Value = 71433830
$ws.Range("I86").Value = 142860900
$ws.Range("J86").Value = 6771.143
$ws.Range("K86").Value = 142860900
$ws.Range("L86").Value = 6771.143
$ws.Range("M86").Value = -142859777
$ws.Range("N86").Value = -9017.143
$ws.Range("H89").Value = 71433830
$ws.Range("I89").Value = 142860900
$ws.Range("J89").Value = 6771.143
$ws.Range("K89").Value = 714304500
$ws.Range("L89").Value = 33855.715
$ws.Range("M89").Value = -714298884
$ws.Range("N89").Value = -45087.715
$ws.Range("H99").Value = 3734.5625
$ws.Range("I99").Value = 1458.4
$ws.Range("J99").Value = 4769.1816
$ws.Range("K99").Value = 1458.4
$ws.Range("L99").Value = 4769.1816
$ws.Range("M99").Value = 39.59999999999991
$ws.Range("N99").Value = -7765.1816
$ws.Range("H126").Value = 3734.5625
$ws.Range("I126").Value = 1458.4
$ws.Range("J126").Value = 4769.1816
$ws.Range("K126").Value = 4375.200000000001
$ws.Range("L126").Value = 14307.5448
$ws.Range("M126").Value = -1905.200000000001
$ws.Range("N126").Value = -19247.5448
$ws.Range("H132").Value = 1694.1765
$ws.Range("I132").Value = 1694.1765
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5082.529500000001
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2552.529500000001
$ws.Range("H134").Value = 1742.6666
$ws.Range("I134").Value = 1610.3158
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 4830.9474
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -2295.9474
$ws.Range("N134").Value = -14070

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1899.56
$ws.Range("I131").Value = 1322.625
$ws.Range("J131").Value = 2009.4524
$ws.Range("K131").Value = 3967.875
$ws.Range("L131").Value = 6028.357199999999
$ws.Range("M131").Value = 1072.125
$ws.Range("N131").Value = -16108.3572

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 72615.11
$ws.Range("I70").Value = 143389.25
$ws.Range("J70").Value = 15995.8
$ws.Range("K70").Value = 143389.25
$ws.Range("L70").Value = 15995.8
$ws.Range("M70").Value = -143119.25
$ws.Range("N70").Value = -16535.8
$ws.Range("H73").Value = 72615.11
$ws.Range("I73").Value = 143389.25
$ws.Range("J73").Value = 15995.8
$ws.Range("K73").Value = 143389.25
$ws.Range("L73").Value = 15995.8
$ws.Range("M73").Value = -142453.25
$ws.Range("N73").Value = -17867.8
$ws.Range("H80").Value = 3635.9285
$ws.Range("I80").Value = 3249.5
$ws.Range("J80").Value = 3925.75
$ws.Range("K80").Value = 3249.5
$ws.Range("L80").Value = 3925.75
$ws.Range("M80").Value = -2251.5
$ws.Range("N80").Value = -5921.75
$ws.Range("H83").Value = 3635.9285
$ws.Range("I83").Value = 3249.5
$ws.Range("J83").Value = 3925.75
$ws.Range("K83").Value = 16247.5
$ws.Range("L83").Value = 19628.75
$ws.Range("M83").Value = -11255.5
$ws.Range("N83").Value = -29612.75
$ws.Range("H132").Value = 4092.1428
$ws.Range("I132").Value = 3764.7827
$ws.Range("K132").Value = 11294.3481
$ws.Range("M132").Value = -8764.348100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4718.9585
$ws.Range("I7").Value = 2605.9167
$ws.Range("K7").Value = 2605.9167
$ws.Range("M7").Value = -2493.9167
$ws.Range("H46").Value = 2014.1666
$ws.Range("I46").Value = 1396.5
$ws.Range("K46").Value = 1396.5
$ws.Range("M46").Value = -1208.5
$ws.Range("H82").Value = 1910.5264
$ws.Range("I82").Value = 1128.0834
$ws.Range("J82").Value = 3251.8572
$ws.Range("K82").Value = 1128.0834
$ws.Range("L82").Value = 3251.8572
$ws.Range("M82").Value = -767.0834
$ws.Range("N82").Value = -3973.8572
$ws.Range("H85").Value = 1910.5264
$ws.Range("I85").Value = 1128.0834
$ws.Range("J85").Value = 3251.8572
$ws.Range("K85").Value = 1128.0834
$ws.Range("L85").Value = 3251.8572
$ws.Range("M85").Value = 119.9166
$ws.Range("N85").Value = -5747.8572
$ws.Range("H93").Value = 3903.9656
$ws.Range("I93").Value = 3441.182
$ws.Range("J93").Value = 5358.4287
$ws.Range("K93").Value = 3441.182
$ws.Range("L93").Value = 5358.4287
$ws.Range("M93").Value = -2193.182
$ws.Range("N93").Value = -7854.4287
$ws.Range("H126").Value = 4718.9585
$ws.Range("I126").Value = 2605.9167
$ws.Range("K126").Value = 7817.750100000001
$ws.Range("M126").Value = -5347.750100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6946.75
$ws.Range("I81").Value = 3655.4285
$ws.Range("J81").Value = 11554.6
$ws.Range("K81").Value = 7310.857
$ws.Range("L81").Value = 23109.2
$ws.Range("M81").Value = -6249.857
$ws.Range("N81").Value = -25231.2
$ws.Range("H84").Value = 6946.75
$ws.Range("I84").Value = 3655.4285
$ws.Range("J84").Value = 11554.6
$ws.Range("K84").Value = 36554.285
$ws.Range("L84").Value = 115546
$ws.Range("M84").Value = -31250.285
$ws.Range("N84").Value = -126154
